# Auto-generated script applying the Odin_Profits.xlsx market-data refresh
# across all 8 item-category worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1395.0435
$ws.Range("I9").Value = 1230.9445
$ws.Range("K9").Value = 1230.9445
$ws.Range("M9").Value = -1061.9445
$ws.Range("H12").Value = 861.4375
$ws.Range("J12").Value = 2420.1667
$ws.Range("L12").Value = 2420.1667
$ws.Range("N12").Value = -2760.1667
$ws.Range("H15").Value = 107942.766
$ws.Range("I15").Value = 107942.766
$ws.Range("K15").Value = 323828.298
$ws.Range("M15").Value = -323659.298
$ws.Range("H55").Value = 199.61539
$ws.Range("I55").Value = 180.42857
$ws.Range("J55").Value = 222
$ws.Range("K55").Value = 180.42857
$ws.Range("L55").Value = 222
$ws.Range("M55").Value = 33.57142999999999
$ws.Range("N55").Value = -650
$ws.Range("H69").Value = 14633.091
$ws.Range("I69").Value = 13441.667
$ws.Range("J69").Value = 19994.5
$ws.Range("K69").Value = 40325.001
$ws.Range("L69").Value = 59983.5
$ws.Range("M69").Value = -39451.001
$ws.Range("N69").Value = -61731.5
$ws.Range("H72").Value = 14633.091
$ws.Range("I72").Value = 13441.667
$ws.Range("J72").Value = 19994.5
$ws.Range("K72").Value = 120975.003
$ws.Range("L72").Value = 179950.5
$ws.Range("M72").Value = -116607.003
$ws.Range("N72").Value = -188686.5
$ws.Range("H116").Value = 12350067
$ws.Range("J116").Value = 5367.1665
$ws.Range("L116").Value = 5367.1665
$ws.Range("N116").Value = -12251.1665
$ws.Range("H132").Value = 560050.7
$ws.Range("I132").Value = 688777.1
$ws.Range("J132").Value = 19399.6
$ws.Range("K132").Value = 2066331.3
$ws.Range("L132").Value = 58198.8
$ws.Range("M132").Value = -2063801.3
$ws.Range("N132").Value = -63258.8
$ws.Range("H138").Value = 4187.048
$ws.Range("J138").Value = 4232.2266
$ws.Range("L138").Value = 12696.6798
$ws.Range("N138").Value = -22976.6798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2517.476
$ws.Range("I45").Value = 2189.2354
$ws.Range("J45").Value = 3912.5
$ws.Range("K45").Value = 2189.2354
$ws.Range("L45").Value = 3912.5
$ws.Range("M45").Value = -1812.2354
$ws.Range("N45").Value = -4666.5
$ws.Range("H61").Value = 4130.1885
$ws.Range("I61").Value = 3723.5312
$ws.Range("J61").Value = 4749.857
$ws.Range("K61").Value = 3723.5312
$ws.Range("L61").Value = 4749.857
$ws.Range("M61").Value = -3511.5312
$ws.Range("N61").Value = -5173.857
$ws.Range("H132").Value = 675831.4399999999
$ws.Range("J132").Value = 40248.5
$ws.Range("L132").Value = 120745.5
$ws.Range("N132").Value = -125805.5
$ws.Range("H136").Value = 4130.1885
$ws.Range("I136").Value = 3723.5312
$ws.Range("J136").Value = 4749.857
$ws.Range("K136").Value = 11170.5936
$ws.Range("L136").Value = 14249.571
$ws.Range("M136").Value = -8620.5936
$ws.Range("N136").Value = -19349.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 14504509
$ws.Range("I80").Value = 1412.091
$ws.Range("K80").Value = 1412.091
$ws.Range("M80").Value = -414.0909999999999
$ws.Range("H83").Value = 14504509
$ws.Range("I83").Value = 1412.091
$ws.Range("K83").Value = 7060.455
$ws.Range("M83").Value = -2068.455
$ws.Range("H105").Value = 3122.9092
$ws.Range("I105").Value = 3206.3044
$ws.Range("K105").Value = 3206.3044
$ws.Range("M105").Value = -1459.3044
$ws.Range("H134").Value = 1284888.5
$ws.Range("I134").Value = 1373414.9
$ws.Range("K134").Value = 4120244.7
$ws.Range("M134").Value = -4117709.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8128.394
$ws.Range("J31").Value = 8975.173000000001
$ws.Range("L31").Value = 8975.173000000001
$ws.Range("N31").Value = -9565.173000000001
$ws.Range("H34").Value = 8128.394
$ws.Range("J34").Value = 8975.173000000001
$ws.Range("L34").Value = 8975.173000000001
$ws.Range("N34").Value = -9379.173000000001
$ws.Range("H58").Value = 5749.143
$ws.Range("J58").Value = 8588
$ws.Range("L58").Value = 8588
$ws.Range("N58").Value = -8994
$ws.Range("H94").Value = 3359.05
$ws.Range("I94").Value = 1031.8334
$ws.Range("J94").Value = 4356.4287
$ws.Range("K94").Value = 1031.8334
$ws.Range("L94").Value = 4356.4287
$ws.Range("M94").Value = -580.8334
$ws.Range("N94").Value = -5258.4287
$ws.Range("H105").Value = 58826684
$ws.Range("I105").Value = 62503292
$ws.Range("K105").Value = 62503292
$ws.Range("M105").Value = -62501545
$ws.Range("H107").Value = 1141
$ws.Range("I107").Value = 1091.5714
$ws.Range("K107").Value = 1091.5714
$ws.Range("M107").Value = 828.4286
$ws.Range("H136").Value = 5749.143
$ws.Range("J136").Value = 8588
$ws.Range("L136").Value = 25764
$ws.Range("N136").Value = -30864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1462.8857
$ws.Range("I5").Value = 844.3333
$ws.Range("J5").Value = 1926.8
$ws.Range("K5").Value = 2532.9999
$ws.Range("L5").Value = 5780.4
$ws.Range("M5").Value = -2420.9999
$ws.Range("N5").Value = -6004.4
$ws.Range("H131").Value = 3271.4211
$ws.Range("I131").Value = 4253.9375
$ws.Range("J131").Value = 2556.8635
$ws.Range("K131").Value = 12761.8125
$ws.Range("L131").Value = 7670.5905
$ws.Range("M131").Value = -7721.8125
$ws.Range("N131").Value = -17750.5905
$ws.Range("H135").Value = 1462.8857
$ws.Range("I135").Value = 844.3333
$ws.Range("J135").Value = 1926.8
$ws.Range("K135").Value = 7598.9997
$ws.Range("L135").Value = 17341.2
$ws.Range("M135").Value = -5063.9997
$ws.Range("N135").Value = -22411.2
$ws.Range("H137").Value = 3071.625
$ws.Range("I137").Value = 2946.4285
$ws.Range("K137").Value = 8839.2855
$ws.Range("M137").Value = -3739.2855
$ws.Range("H139").Value = 31252742
$ws.Range("I139").Value = 50001496
$ws.Range("K139").Value = 150004488
$ws.Range("M139").Value = -149999348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 9500
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 9000
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -9755
$ws.Range("N20").Value = -9490
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H24").Value = 2502624.8
$ws.Range("I24").Value = 2502624.8
$ws.Range("K24").Value = 2502624.8
$ws.Range("M24").Value = -2502451.8
$ws.Range("H25").Value = 25000
$ws.Range("J25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("N25").Value = -26058
$ws.Range("H70").Value = 5924.8
$ws.Range("I70").Value = 5732.5454
$ws.Range("J70").Value = 6075.857
$ws.Range("K70").Value = 5732.5454
$ws.Range("L70").Value = 6075.857
$ws.Range("M70").Value = -5462.5454
$ws.Range("N70").Value = -6615.857
$ws.Range("H73").Value = 5924.8
$ws.Range("I73").Value = 5732.5454
$ws.Range("J73").Value = 6075.857
$ws.Range("K73").Value = 5732.5454
$ws.Range("L73").Value = 6075.857
$ws.Range("M73").Value = -4796.5454
$ws.Range("N73").Value = -7947.857
$ws.Range("H137").Value = 100747.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 199999.33
$ws.Range("J139").Value = 199999.33
$ws.Range("L139").Value = 199999.33
$ws.Range("N139").Value = -210279.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1753267.4
$ws.Range("J43").Value = 1992230.8
$ws.Range("L43").Value = 1992230.8
$ws.Range("N43").Value = -1992616.8
$ws.Range("H46").Value = 33335696
$ws.Range("I46").Value = 872.5
$ws.Range("K46").Value = 872.5
$ws.Range("M46").Value = -684.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 4546536.5
$ws.Range("I23").Value = 5001100
$ws.Range("J23").Value = 899
$ws.Range("K23").Value = 5001100
$ws.Range("L23").Value = 899
$ws.Range("M23").Value = -5000871
$ws.Range("N23").Value = -1357
$ws.Range("H104").Value = 144997
$ws.Range("J104").Value = 144997
$ws.Range("L104").Value = 144997
$ws.Range("N104").Value = -151985
$ws.Range("H126").Value = 6446.4375
$ws.Range("J126").Value = 8262.916999999999
$ws.Range("L126").Value = 24788.751
$ws.Range("N126").Value = -29728.751

